# Apply updated "想去人数" (want-to-go count) figures to the F column
# across the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types)
# sheets, matching the published gh-pages regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 111
$ws1.Range("F3").Value  = 12219
$ws1.Range("F4").Value  = 55
$ws1.Range("F5").Value  = 239
$ws1.Range("F6").Value  = 376
$ws1.Range("F8").Value  = 12146
$ws1.Range("F9").Value  = 510
$ws1.Range("F10").Value = 1194
$ws1.Range("F12").Value = 611
$ws1.Range("F13").Value = 2808
$ws1.Range("F14").Value = 5966
$ws1.Range("F16").Value = 3570

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 13

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 111
$ws4.Range("F4").Value  = 12219
$ws4.Range("F5").Value  = 55
$ws4.Range("F6").Value  = 239
$ws4.Range("F7").Value  = 13
$ws4.Range("F8").Value  = 376
$ws4.Range("F10").Value = 12147
$ws4.Range("F11").Value = 510
$ws4.Range("F12").Value = 1194
$ws4.Range("F14").Value = 612
$ws4.Range("F15").Value = 2808
$ws4.Range("F17").Value = 5966
$ws4.Range("F19").Value = 3570
